$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dragon Riches, an online slot game with Chinese lucky charms and five golden symbols. Play for free and experience the oriental theme.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document: remove the duplicated bold title
#    paragraph ("Play Dragon Riches Free: ...") that used to sit just
#    before the meta-description paragraph.
# ---------------------------------------------------------------------------
$target = "Play Dragon Riches Free: Chinese Themed Slot Game Review"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq $target) {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph - what used to be the
#    meta-description text is now repurposed as an AI image-generation
#    prompt.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = "Prompt: Create a cartoon-style image to feature in a blog post about the online slot game ""Dragon Riches"". The image should prominently feature a happy Maya warrior with glasses. The warrior should be holding a golden dragon and surrounded by Chinese lucky charms, such as coins, paper lanterns, and yuanbao ingots. The background should include elements of both Mayan and Chinese culture, such as temples and dragons. The image should give off a cheerful, lucky, and adventurous vibe to entice readers to give the game a try."
